# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H, I, J, K, L, M, N) across the Leve profit sheets with freshly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 573.75
$ws.Range("J2").Value = 633.3333
$ws.Range("L2").Value = 633.3333
$ws.Range("N2").Value = -859.3333

$ws.Range("H33").Value = 43478664
$ws.Range("I33").Value = 47619400
$ws.Range("K33").Value = 47619400
$ws.Range("M33").Value = -47619171

$ws.Range("H98").Value = 1153.3043
$ws.Range("I98").Value = 1249.6666
$ws.Range("J98").Value = 1048.1818
$ws.Range("K98").Value = 1249.6666
$ws.Range("L98").Value = 1048.1818
$ws.Range("M98").Value = 248.3334
$ws.Range("N98").Value = -4044.1818

$ws.Range("H113").Value = 216747.28
$ws.Range("I113").Value = 364307.75
$ws.Range("J113").Value = 20000
$ws.Range("K113").Value = 364307.75
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = -361053.75
$ws.Range("N113").Value = -26508

$ws.Range("H122").Value = 1153.3043
$ws.Range("I122").Value = 1249.6666
$ws.Range("J122").Value = 1048.1818
$ws.Range("K122").Value = 3748.9998
$ws.Range("L122").Value = 3144.5454
$ws.Range("M122").Value = -1298.9998
$ws.Range("N122").Value = -8044.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10000
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = -9885
$ws.Range("N3").Value = 0

$ws.Range("H110").Value = 1344.4
$ws.Range("I110").Value = 1168.4546
$ws.Range("K110").Value = 1168.4546
$ws.Range("M110").Value = 876.5454

$ws.Range("H122").Value = 10797.333
$ws.Range("I122").Value = 23042.4
$ws.Range("K122").Value = 69127.20000000001
$ws.Range("M122").Value = -66677.20000000001

$ws.Range("H137").Value = 5884808.5
$ws.Range("J137").Value = 15388268
$ws.Range("L137").Value = 46164804
$ws.Range("N137").Value = -46169904

$ws.Range("H138").Value = 5210712
$ws.Range("I138").Value = 1279.0571
$ws.Range("K138").Value = 3837.1713
$ws.Range("M138").Value = 1302.8287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 764025.4
$ws.Range("I7").Value = 858886.1
$ws.Range("J7").Value = 100000
$ws.Range("K7").Value = 858886.1
$ws.Range("L7").Value = 100000
$ws.Range("M7").Value = -858773.1
$ws.Range("N7").Value = -100226

$ws.Range("H10").Value = 25102.5
$ws.Range("I10").Value = 205
$ws.Range("J10").Value = 50000
$ws.Range("K10").Value = 205
$ws.Range("L10").Value = 50000
$ws.Range("M10").Value = -65
$ws.Range("N10").Value = -50280

$ws.Range("H86").Value = 19232538
$ws.Range("I86").Value = 1803.5385
$ws.Range("J86").Value = 38463270
$ws.Range("K86").Value = 1803.5385
$ws.Range("L86").Value = 38463270
$ws.Range("M86").Value = -680.5385000000001
$ws.Range("N86").Value = -38465516

$ws.Range("H89").Value = 19232538
$ws.Range("I89").Value = 1803.5385
$ws.Range("J89").Value = 38463270
$ws.Range("K89").Value = 9017.692500000001
$ws.Range("L89").Value = 192316350
$ws.Range("M89").Value = -3401.692500000001
$ws.Range("N89").Value = -192327582

$ws.Range("H94").Value = 1053.6111
$ws.Range("I94").Value = 595
$ws.Range("J94").Value = 1774.2858
$ws.Range("K94").Value = 595
$ws.Range("L94").Value = 1774.2858
$ws.Range("M94").Value = -144
$ws.Range("N94").Value = -2676.2858

$ws.Range("H107").Value = 1794.4193
$ws.Range("I107").Value = 1941.3636
$ws.Range("J107").Value = 1435.2222
$ws.Range("K107").Value = 1941.3636
$ws.Range("L107").Value = 1435.2222
$ws.Range("M107").Value = -21.36359999999991
$ws.Range("N107").Value = -5275.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 33802
$ws.Range("J70").Value = 33802
$ws.Range("L70").Value = 33802
$ws.Range("N70").Value = -34432

$ws.Range("H73").Value = 33802
$ws.Range("J73").Value = 33802
$ws.Range("L73").Value = 33802
$ws.Range("N73").Value = -35986

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2182.4
$ws.Range("I64").Value = 970.6667
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 2912.0001
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -2642.0001
$ws.Range("N64").Value = -12540

$ws.Range("H67").Value = 2182.4
$ws.Range("I67").Value = 970.6667
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 2912.0001
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -1976.0001
$ws.Range("N67").Value = -13872

$ws.Range("H68").Value = 1136.8628
$ws.Range("I68").Value = 867.0417
$ws.Range("J68").Value = 1376.7037
$ws.Range("K68").Value = 2601.1251
$ws.Range("L68").Value = 4130.1111
$ws.Range("M68").Value = -1790.1251
$ws.Range("N68").Value = -5752.1111

$ws.Range("H71").Value = 1136.8628
$ws.Range("I71").Value = 867.0417
$ws.Range("J71").Value = 1376.7037
$ws.Range("K71").Value = 7803.3753
$ws.Range("L71").Value = 12390.3333
$ws.Range("M71").Value = -3747.3753
$ws.Range("N71").Value = -20502.3333

$ws.Range("H107").Value = 938.9245
$ws.Range("I107").Value = 325.51852
$ws.Range("J107").Value = 1575.9231
$ws.Range("K107").Value = 976.55556
$ws.Range("L107").Value = 4727.7693
$ws.Range("M107").Value = 943.44444
$ws.Range("N107").Value = -8567.7693

$ws.Range("H113").Value = 699.6269
$ws.Range("I113").Value = 564.4820999999999
$ws.Range("J113").Value = 1387.6364
$ws.Range("K113").Value = 1693.4463
$ws.Range("L113").Value = 4162.9092
$ws.Range("M113").Value = 476.5537000000002
$ws.Range("N113").Value = -8502.9092

$ws.Range("H122").Value = 1731.1111
$ws.Range("J122").Value = 756.5714
$ws.Range("L122").Value = 6809.1426
$ws.Range("N122").Value = -11709.1426

$ws.Range("H141").Value = 13347.167
$ws.Range("J141").Value = 13347.167
$ws.Range("L141").Value = 13347.167
$ws.Range("N141").Value = -23707.167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 9
$ws.Range("J6").Value = 9
$ws.Range("L6").Value = 9
$ws.Range("N6").Value = -235

$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 200
$ws.Range("K7").Value = 200
$ws.Range("M7").Value = -88

$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 200
$ws.Range("K8").Value = 200
$ws.Range("M8").Value = -61

$ws.Range("H9").Value = 278.5
$ws.Range("I9").Value = 278.5
$ws.Range("K9").Value = 278.5
$ws.Range("M9").Value = -108.5

$ws.Range("H13").Value = 549.5
$ws.Range("J13").Value = 699.3333
$ws.Range("L13").Value = 699.3333
$ws.Range("N13").Value = -977.3333

$ws.Range("H16").Value = 9
$ws.Range("J16").Value = 9
$ws.Range("L16").Value = 9
$ws.Range("N16").Value = -509

$ws.Range("H19").Value = 6773
$ws.Range("J19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("N19").Value = -10576

$ws.Range("H23").Value = 3353.3333
$ws.Range("J23").Value = 3353.3333
$ws.Range("L23").Value = 3353.3333
$ws.Range("N23").Value = -3799.3333

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3220.9524
$ws.Range("I16").Value = 2575.7144
$ws.Range("J16").Value = 4511.4287
$ws.Range("K16").Value = 2575.7144
$ws.Range("L16").Value = 4511.4287
$ws.Range("M16").Value = -2405.7144
$ws.Range("N16").Value = -4851.4287

$ws.Range("H22").Value = 1008.9231
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1014.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1014.5
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1604.5

$ws.Range("H27").Value = 1008.9231
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1014.5
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1014.5
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1228.5

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = 0

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = 0
